# Show Cam Renewal -- rebuild Sheet1 with the new compact "ExecSet" layout
# and move the live selection on both sheets.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# ---------------------------------------------------------------------------
# 1) Pre-register the brand-new shared strings in the exact order Excel used
#    when this sheet was authored (scratch column far off-sheet, then
#    Clear() it -- the shared-string slot survives because every one of
#    these strings gets reused below, it just keeps its original index).
# ---------------------------------------------------------------------------
$newStringOrder = @(
    "---",
    "bools",
    "FormatVer",
    "DataVer",
    "StartLine",
    "Duration",
    "S:EffectData:Int32",
    "S:EffectData:String",
    "A(S:EffectData:Int32)",
    "A(S:EffectData:String)",
    "EffectData.ID",
    "EffectData.Name",
    "EffectDatas.ID",
    "EffectDatas.Name",
    "ExecSet.EffectData.ID",
    "ExecSet.EffectData.Name",
    "ExecSet.EffectDatas.ID",
    "ExecSet.EffectDatas.Name",
    "ExecSet.ID",
    "S:ExecSet:Int32",
    "S:ExecSet:EffectData:Int32",
    "S:ExecSet:EffectData:String",
    "S:ExecSet:A(S:EffectData:Int32)",
    "S:ExecSet:A(S:EffectData:String)"
)

$scratchRow = 1
foreach ($s in $newStringOrder) {
    $ws1.Cells.Item($scratchRow, 26).Value = $s   # column Z, far away scratch area
    $scratchRow = $scratchRow + 1
}
$ws1.Range($ws1.Cells.Item(1,26), $ws1.Cells.Item($newStringOrder.Count,26)).Clear()

# ---------------------------------------------------------------------------
# 2) Wipe the old Sheet1 content and lay down the new 15 x 2 table.
# ---------------------------------------------------------------------------
$ws1.Cells.Clear()

# Row 4's quote-prefixed cell must claim its cellXfs slot (index 1) before
# the centered numeric cells below claim theirs (index 2) -- matches the
# style-table order the original author ended up with.
$ws1.Cells.Item(4,1).Value = "RowName"
$ws1.Cells.Item(4,2).Value = "'---"

$ws1.Cells.Item(1,1).Value = "FormatVer"
$ws1.Cells.Item(1,2).Value = 1
$ws1.Cells.Item(1,2).HorizontalAlignment = -4108

$ws1.Cells.Item(2,1).Value = "DataVer"
$ws1.Cells.Item(2,2).Value = 1
$ws1.Cells.Item(2,2).HorizontalAlignment = -4108

$ws1.Cells.Item(3,1).Value = "StartLine"
$ws1.Cells.Item(3,2).Value = 1
$ws1.Cells.Item(3,2).HorizontalAlignment = -4108

$ws1.Cells.Item(5,1).Value = "Duration"
$ws1.Cells.Item(5,2).Value = "float"

$ws1.Cells.Item(6,1).Value = "EffectData.ID"
$ws1.Cells.Item(6,2).Value = "S:EffectData:Int32"

$ws1.Cells.Item(7,1).Value = "EffectData.Name"
$ws1.Cells.Item(7,2).Value = "S:EffectData:String"

$ws1.Cells.Item(8,1).Value = "bools"
$ws1.Cells.Item(8,2).Value = "A:V:bool"

$ws1.Cells.Item(9,1).Value = "EffectDatas.ID"
$ws1.Cells.Item(9,2).Value = "A(S:EffectData:Int32)"

$ws1.Cells.Item(10,1).Value = "EffectDatas.Name"
$ws1.Cells.Item(10,2).Value = "A(S:EffectData:String)"

$ws1.Cells.Item(11,1).Value = "ExecSet.ID"
$ws1.Cells.Item(11,2).Value = "S:ExecSet:Int32"

$ws1.Cells.Item(12,1).Value = "ExecSet.EffectData.ID"
$ws1.Cells.Item(12,2).Value = "S:ExecSet:EffectData:Int32"

$ws1.Cells.Item(13,1).Value = "ExecSet.EffectData.Name"
$ws1.Cells.Item(13,2).Value = "S:ExecSet:EffectData:String"

$ws1.Cells.Item(14,1).Value = "ExecSet.EffectDatas.ID"
$ws1.Cells.Item(14,2).Value = "S:ExecSet:A(S:EffectData:Int32)"

$ws1.Cells.Item(15,1).Value = "ExecSet.EffectDatas.Name"
$ws1.Cells.Item(15,2).Value = "S:ExecSet:A(S:EffectData:String)"

# Column widths (author resized A/B to fit the new, longer labels).
$ws1.Columns.Item(1).ColumnWidth = 27.84
$ws1.Columns.Item(2).ColumnWidth = 35.41

# ---------------------------------------------------------------------------
# 3) View/selection changes.
#    Sheet2 keeps its original data untouched; only its live selection moves.
#    Select Sheet2 first so the final active sheet/tab ends up back on Sheet1.
# ---------------------------------------------------------------------------
$ws2.Range("C16").Select()
$ws1.Range("B17:B18").Select()
